$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 41668420
$ws.Range("J18").Value = 2000
$ws.Range("L18").Value = 2000
$ws.Range("N18").Value = -2568
$ws.Range("H38").Value = 447.33334
$ws.Range("I38").Value = 136.8
$ws.Range("K38").Value = 410.4
$ws.Range("M38").Value = -38.40000000000003
$ws.Range("H39").Value = 542.125
$ws.Range("I39").Value = 225
$ws.Range("J39").Value = 1916.3334
$ws.Range("K39").Value = 675
$ws.Range("L39").Value = 5749.0002
$ws.Range("M39").Value = -379
$ws.Range("N39").Value = -6341.0002
$ws.Range("H43").Value = 1856
$ws.Range("I43").Value = 1822
$ws.Range("K43").Value = 1822
$ws.Range("M43").Value = -1753
$ws.Range("H54").Value = 6297
$ws.Range("I54").Value = 5062.6665
$ws.Range("J54").Value = 10000
$ws.Range("K54").Value = 5062.6665
$ws.Range("L54").Value = 10000
$ws.Range("M54").Value = -4576.6665
$ws.Range("N54").Value = -10972
$ws.Range("H55").Value = 286.5
$ws.Range("I55").Value = 106.77778
$ws.Range("J55").Value = 610
$ws.Range("K55").Value = 106.77778
$ws.Range("L55").Value = 610
$ws.Range("M55").Value = 107.22222
$ws.Range("N55").Value = -1038
$ws.Range("H62").Value = 2233
$ws.Range("I62").Value = 1899.5
$ws.Range("K62").Value = 1899.5
$ws.Range("M62").Value = -1275.5
$ws.Range("H64").Value = 3927.5
$ws.Range("I64").Value = 3325
$ws.Range("J64").Value = 5333.3335
$ws.Range("K64").Value = 3325
$ws.Range("L64").Value = 5333.3335
$ws.Range("M64").Value = -3077
$ws.Range("N64").Value = -5829.3335
$ws.Range("H65").Value = 2233
$ws.Range("I65").Value = 1899.5
$ws.Range("K65").Value = 9497.5
$ws.Range("M65").Value = -6377.5
$ws.Range("H67").Value = 3927.5
$ws.Range("I67").Value = 3325
$ws.Range("J67").Value = 5333.3335
$ws.Range("K67").Value = 3325
$ws.Range("L67").Value = 5333.3335
$ws.Range("M67").Value = -2467
$ws.Range("N67").Value = -7049.3335
$ws.Range("H74").Value = 5000
$ws.Range("J74").Value = 5000
$ws.Range("L74").Value = 5000
$ws.Range("N74").Value = -6872
$ws.Range("H76").Value = 4424.5
$ws.Range("I76").Value = 4599
$ws.Range("J76").Value = 4250
$ws.Range("K76").Value = 4599
$ws.Range("L76").Value = 4250
$ws.Range("M76").Value = -4284
$ws.Range("N76").Value = -4880
$ws.Range("H77").Value = 5000
$ws.Range("J77").Value = 5000
$ws.Range("L77").Value = 25000
$ws.Range("N77").Value = -34360
$ws.Range("H79").Value = 4424.5
$ws.Range("I79").Value = 4599
$ws.Range("J79").Value = 4250
$ws.Range("K79").Value = 4599
$ws.Range("L79").Value = 4250
$ws.Range("M79").Value = -3507
$ws.Range("N79").Value = -6434
$ws.Range("H86").Value = 23733.8
$ws.Range("I86").Value = 4037.6667
$ws.Range("K86").Value = 4037.6667
$ws.Range("M86").Value = -2914.6667
$ws.Range("H89").Value = 23733.8
$ws.Range("I89").Value = 4037.6667
$ws.Range("K89").Value = 20188.3335
$ws.Range("M89").Value = -14572.3335
$ws.Range("H101").Value = 466
$ws.Range("I101").Value = 199
$ws.Range("K101").Value = 597
$ws.Range("M101").Value = 1025
$ws.Range("H132").Value = 8393.73
$ws.Range("I132").Value = 5522.875
$ws.Range("K132").Value = 16568.625
$ws.Range("M132").Value = -14038.625
$ws.Range("H137").Value = 1168.3704
$ws.Range("I137").Value = 1389.875
$ws.Range("J137").Value = 846.1818
$ws.Range("K137").Value = 4169.625
$ws.Range("L137").Value = 2538.5454
$ws.Range("M137").Value = -1619.625
$ws.Range("N137").Value = -7638.5454
$ws.Range("H140").Value = 42093
$ws.Range("I140").Value = 44549
$ws.Range("J140").Value = 41547.223
$ws.Range("K140").Value = 44549
$ws.Range("L140").Value = 41547.223
$ws.Range("M140").Value = -39369
$ws.Range("N140").Value = -51907.223

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2673.4387
$ws.Range("I32").Value = 1669.6522
$ws.Range("J32").Value = 18064.834
$ws.Range("K32").Value = 1669.6522
$ws.Range("L32").Value = 18064.834
$ws.Range("M32").Value = -1382.6522
$ws.Range("N32").Value = -18638.834
$ws.Range("H45").Value = 8206
$ws.Range("I45").Value = 11591.833
$ws.Range("J45").Value = 3691.5557
$ws.Range("K45").Value = 11591.833
$ws.Range("L45").Value = 3691.5557
$ws.Range("M45").Value = -11214.833
$ws.Range("N45").Value = -4445.5557
$ws.Range("H74").Value = 5668.5
$ws.Range("I74").Value = 5914.4414
$ws.Range("J74").Value = 4623.25
$ws.Range("K74").Value = 5914.4414
$ws.Range("L74").Value = 4623.25
$ws.Range("M74").Value = -5040.4414
$ws.Range("N74").Value = -6371.25
$ws.Range("H77").Value = 5668.5
$ws.Range("I77").Value = 5914.4414
$ws.Range("J77").Value = 4623.25
$ws.Range("K77").Value = 29572.207
$ws.Range("L77").Value = 23116.25
$ws.Range("M77").Value = -25204.207
$ws.Range("N77").Value = -31852.25
$ws.Range("H88").Value = 983.3333
$ws.Range("I88").Value = 933.3333
$ws.Range("K88").Value = 933.3333
$ws.Range("M88").Value = -527.3333
$ws.Range("H91").Value = 983.3333
$ws.Range("I91").Value = 933.3333
$ws.Range("K91").Value = 933.3333
$ws.Range("M91").Value = 470.6667
$ws.Range("H93").Value = 34995
$ws.Range("J93").Value = 34995
$ws.Range("L93").Value = 34995
$ws.Range("N93").Value = -39987
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
$ws.Range("H111").Value = 89319.664
$ws.Range("J111").Value = 89319.664
$ws.Range("L111").Value = 89319.664
$ws.Range("N111").Value = -97499.664
$ws.Range("H122").Value = 2746.2778
$ws.Range("I122").Value = 2769.75
$ws.Range("J122").Value = 2699.3333
$ws.Range("K122").Value = 8309.25
$ws.Range("L122").Value = 8097.999899999999
$ws.Range("M122").Value = -5859.25
$ws.Range("N122").Value = -12997.9999
$ws.Range("H132").Value = 2831.7
$ws.Range("I132").Value = 2676.8572
$ws.Range("K132").Value = 8030.571599999999
$ws.Range("M132").Value = -5500.571599999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H6").Value = 99990
$ws.Range("J6").Value = 99990
$ws.Range("L6").Value = 99990
$ws.Range("N6").Value = -100216
$ws.Range("H86").Value = 2449.1538
$ws.Range("I86").Value = 2167.6365
$ws.Range("J86").Value = 3997.5
$ws.Range("K86").Value = 2167.6365
$ws.Range("L86").Value = 3997.5
$ws.Range("M86").Value = -1044.6365
$ws.Range("N86").Value = -6243.5
$ws.Range("H89").Value = 2449.1538
$ws.Range("I89").Value = 2167.6365
$ws.Range("J89").Value = 3997.5
$ws.Range("K89").Value = 10838.1825
$ws.Range("L89").Value = 19987.5
$ws.Range("M89").Value = -5222.182500000001
$ws.Range("N89").Value = -31219.5
$ws.Range("H107").Value = 4186.8813
$ws.Range("I107").Value = 3919.1428
$ws.Range("J107").Value = 4428.7095
$ws.Range("K107").Value = 3919.1428
$ws.Range("L107").Value = 4428.7095
$ws.Range("M107").Value = -1999.1428
$ws.Range("N107").Value = -8268.709500000001
$ws.Range("H114").Value = 59999
$ws.Range("J114").Value = 59999
$ws.Range("L114").Value = 59999
$ws.Range("N114").Value = -68677
$ws.Range("H133").Value = 86523.5
$ws.Range("J133").Value = 85338
$ws.Range("L133").Value = 85338
$ws.Range("N133").Value = -95458
$ws.Range("H134").Value = 4498.943
$ws.Range("I134").Value = 4202
$ws.Range("K134").Value = 12606
$ws.Range("M134").Value = -10071

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 66.57143000000001
$ws.Range("I7").Value = 30.333334
$ws.Range("J7").Value = 93.75
$ws.Range("K7").Value = 30.333334
$ws.Range("L7").Value = 93.75
$ws.Range("M7").Value = 82.66666599999999
$ws.Range("N7").Value = -319.75
$ws.Range("H9").Value = 99961
$ws.Range("J9").Value = 99961
$ws.Range("L9").Value = 99961
$ws.Range("N9").Value = -100297
$ws.Range("H22").Value = 1445.4584
$ws.Range("I22").Value = 1369.6154
$ws.Range("J22").Value = 1535.091
$ws.Range("K22").Value = 1369.6154
$ws.Range("L22").Value = 1535.091
$ws.Range("M22").Value = -1019.6154
$ws.Range("N22").Value = -2235.091
$ws.Range("H31").Value = 2321.95
$ws.Range("I31").Value = 2262.2
$ws.Range("J31").Value = 2381.7
$ws.Range("K31").Value = 2262.2
$ws.Range("L31").Value = 2381.7
$ws.Range("M31").Value = -1967.2
$ws.Range("N31").Value = -2971.7
$ws.Range("H34").Value = 2321.95
$ws.Range("I34").Value = 2262.2
$ws.Range("J34").Value = 2381.7
$ws.Range("K34").Value = 2262.2
$ws.Range("L34").Value = 2381.7
$ws.Range("M34").Value = -2060.2
$ws.Range("N34").Value = -2785.7
$ws.Range("H55").Value = 19900
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 19900
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 19900
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -20530
$ws.Range("H64").Value = 59999
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 59999
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 59999
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -60495
$ws.Range("H67").Value = 59999
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 59999
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 59999
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -61715
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").ClearContents()
$ws.Range("H105").Value = 2031.7142
$ws.Range("I105").Value = 1953.6666
$ws.Range("K105").Value = 1953.6666
$ws.Range("M105").Value = -206.6666
$ws.Range("H107").Value = 1342
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("H134").Value = 7045.184
$ws.Range("I134").Value = 7727.241
$ws.Range("K134").Value = 23181.723
$ws.Range("M134").Value = -20646.723
$ws.Range("H140").Value = 76426.63
$ws.Range("I140").Value = 1709
$ws.Range("J140").Value = 83898.39999999999
$ws.Range("K140").Value = 1709
$ws.Range("L140").Value = 83898.39999999999
$ws.Range("M140").Value = 3471
$ws.Range("N140").Value = -94258.39999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 375598.75
$ws.Range("I9").Value = 460538.9
$ws.Range("J9").Value = 1862
$ws.Range("K9").Value = 1381616.7
$ws.Range("L9").Value = 5586
$ws.Range("M9").Value = -1381392.7
$ws.Range("N9").Value = -6034
$ws.Range("H23").Value = 2345.4167
$ws.Range("I23").Value = 2740.5
$ws.Range("J23").Value = 1950.3334
$ws.Range("K23").Value = 8221.5
$ws.Range("L23").Value = 5851.0002
$ws.Range("M23").Value = -7986.5
$ws.Range("N23").Value = -6321.0002
$ws.Range("H44").Value = 647.6667
$ws.Range("I44").Value = 647.6667
$ws.Range("K44").Value = 1943.0001
$ws.Range("M44").Value = -1545.0001
$ws.Range("H51").Value = 1619.6
$ws.Range("I51").Value = 1499.5
$ws.Range("K51").Value = 4498.5
$ws.Range("M51").Value = -4038.5
$ws.Range("H82").Value = 11160.625
$ws.Range("I82").Value = 8214.5
$ws.Range("J82").Value = 19999
$ws.Range("K82").Value = 24643.5
$ws.Range("L82").Value = 59997
$ws.Range("M82").Value = -24237.5
$ws.Range("N82").Value = -60809
$ws.Range("H85").Value = 11160.625
$ws.Range("I85").Value = 8214.5
$ws.Range("J85").Value = 19999
$ws.Range("K85").Value = 24643.5
$ws.Range("L85").Value = 59997
$ws.Range("M85").Value = -23239.5
$ws.Range("N85").Value = -62805
$ws.Range("H98").Value = 270.15384
$ws.Range("I98").Value = 192.33333
$ws.Range("K98").Value = 576.99999
$ws.Range("M98").Value = 921.00001
$ws.Range("H107").Value = 353.75
$ws.Range("I107").Value = 265
$ws.Range("K107").Value = 795
$ws.Range("M107").Value = 1125
$ws.Range("H113").Value = 546.9375
$ws.Range("I113").Value = 405.14285
$ws.Range("J113").Value = 657.2222
$ws.Range("K113").Value = 1215.42855
$ws.Range("L113").Value = 1971.6666
$ws.Range("M113").Value = 954.5714499999999
$ws.Range("N113").Value = -6311.6666
$ws.Range("H118").Value = 601.8
$ws.Range("I118").Value = 601.8
$ws.Range("K118").Value = 1805.4
$ws.Range("M118").Value = -562.3999999999999
$ws.Range("H119").Value = 6369.857
$ws.Range("I119").Value = 2945
$ws.Range("K119").Value = 8835
$ws.Range("M119").Value = -3997
$ws.Range("H122").Value = 1539.6
$ws.Range("I122").Value = 679.4
$ws.Range("J122").Value = 2399.8
$ws.Range("K122").Value = 6114.599999999999
$ws.Range("L122").Value = 21598.2
$ws.Range("M122").Value = -3664.599999999999
$ws.Range("N122").Value = -26498.2
$ws.Range("H128").Value = 181705.58
$ws.Range("I128").Value = 181705.58
$ws.Range("K128").Value = 545116.74
$ws.Range("M128").Value = -540136.74
$ws.Range("H137").Value = 7758.6875
$ws.Range("I137").Value = 4640.769
$ws.Range("K137").Value = 13922.307
$ws.Range("M137").Value = -8822.307000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6792.4375
$ws.Range("J70").Value = 6436.5
$ws.Range("L70").Value = 6436.5
$ws.Range("N70").Value = -6976.5
$ws.Range("H73").Value = 6792.4375
$ws.Range("J73").Value = 6436.5
$ws.Range("L73").Value = 6436.5
$ws.Range("N73").Value = -8308.5
$ws.Range("H80").Value = 21518740
$ws.Range("I80").Value = 50715860
$ws.Range("J80").Value = 5071.9473
$ws.Range("K80").Value = 50715860
$ws.Range("L80").Value = 5071.9473
$ws.Range("M80").Value = -50714862
$ws.Range("N80").Value = -7067.9473
$ws.Range("H83").Value = 21518740
$ws.Range("I83").Value = 50715860
$ws.Range("J83").Value = 5071.9473
$ws.Range("K83").Value = 253579300
$ws.Range("L83").Value = 25359.7365
$ws.Range("M83").Value = -253574308
$ws.Range("N83").Value = -35343.7365
$ws.Range("H99").Value = 15131.571
$ws.Range("I99").Value = 10731.5
$ws.Range("J99").Value = 20998.334
$ws.Range("K99").Value = 10731.5
$ws.Range("L99").Value = 20998.334
$ws.Range("M99").Value = -8485.5
$ws.Range("N99").Value = -25490.334
$ws.Range("H102").Value = 2991.647
$ws.Range("I102").Value = 3183.25
$ws.Range("K102").Value = 3183.25
$ws.Range("M102").Value = -1561.25
$ws.Range("H113").Value = 3261.4167
$ws.Range("I113").Value = 2925.625
$ws.Range("K113").Value = 2925.625
$ws.Range("M113").Value = -755.625
$ws.Range("H122").Value = 4410.636
$ws.Range("I122").Value = 4545.2856
$ws.Range("J122").Value = 4175
$ws.Range("K122").Value = 13635.8568
$ws.Range("L122").Value = 12525
$ws.Range("M122").Value = -11185.8568
$ws.Range("N122").Value = -17425
$ws.Range("H126").Value = 5414.7026
$ws.Range("I126").Value = 9371.5
$ws.Range("K126").Value = 28114.5
$ws.Range("M126").Value = -25644.5
$ws.Range("H132").Value = 6456.706
$ws.Range("I132").Value = 6554.5713
$ws.Range("J132").Value = 6000
$ws.Range("K132").Value = 19663.7139
$ws.Range("L132").Value = 18000
$ws.Range("M132").Value = -17133.7139
$ws.Range("N132").Value = -23060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 1500
$ws.Range("J17").Value = 1500
$ws.Range("L17").Value = 1500
$ws.Range("N17").Value = -1840
$ws.Range("H40").Value = 18356.445
$ws.Range("I40").Value = 18994.533
$ws.Range("K40").Value = 18994.533
$ws.Range("M40").Value = -18858.533
$ws.Range("H46").Value = 4114.6313
$ws.Range("I46").Value = 3312.375
$ws.Range("J46").Value = 4698.091
$ws.Range("K46").Value = 3312.375
$ws.Range("L46").Value = 4698.091
$ws.Range("M46").Value = -3124.375
$ws.Range("N46").Value = -5074.091
$ws.Range("H68").Value = 15155370
$ws.Range("I68").Value = 18520152
$ws.Range("J68").Value = 13848.75
$ws.Range("K68").Value = 18520152
$ws.Range("L68").Value = 13848.75
$ws.Range("M68").Value = -18519403
$ws.Range("N68").Value = -15346.75
$ws.Range("H71").Value = 15155370
$ws.Range("I71").Value = 18520152
$ws.Range("J71").Value = 13848.75
$ws.Range("K71").Value = 92600760
$ws.Range("L71").Value = 69243.75
$ws.Range("M71").Value = -92597016
$ws.Range("N71").Value = -76731.75
$ws.Range("H93").Value = 16668127
$ws.Range("I93").Value = 18183258
$ws.Range("J93").Value = 1700
$ws.Range("K93").Value = 18183258
$ws.Range("L93").Value = 1700
$ws.Range("M93").Value = -18182010
$ws.Range("N93").Value = -4196
$ws.Range("H100").Value = 50001816
$ws.Range("I100").Value = 76924590
$ws.Range("J100").Value = 2369.8572
$ws.Range("K100").Value = 76924590
$ws.Range("L100").Value = 2369.8572
$ws.Range("M100").Value = -76924049
$ws.Range("N100").Value = -3451.8572
$ws.Range("H116").Value = 97495
$ws.Range("J116").Value = 97495
$ws.Range("L116").Value = 97495
$ws.Range("N116").Value = -106673
$ws.Range("H122").Value = 10820.389
$ws.Range("I122").Value = 11097.714
$ws.Range("J122").Value = 9849.75
$ws.Range("K122").Value = 33293.142
$ws.Range("L122").Value = 29549.25
$ws.Range("M122").Value = -30843.142
$ws.Range("N122").Value = -34449.25
$ws.Range("H132").Value = 55216.176
$ws.Range("I132").Value = 57044.184
$ws.Range("J132").Value = 15000
$ws.Range("K132").Value = 171132.552
$ws.Range("L132").Value = 45000
$ws.Range("M132").Value = -168602.552
$ws.Range("N132").Value = -50060
$ws.Range("H136").Value = 8258.825999999999
$ws.Range("I136").Value = 6586.143
$ws.Range("J136").Value = 10860.777
$ws.Range("K136").Value = 19758.429
$ws.Range("L136").Value = 32582.331
$ws.Range("M136").Value = -17208.429
$ws.Range("N136").Value = -37682.331
$ws.Range("H137").Value = 199999
$ws.Range("J137").Value = 199999
$ws.Range("L137").Value = 199999
$ws.Range("N137").Value = -210199

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 8122.5
$ws.Range("J5").Value = 4947
$ws.Range("L5").Value = 4947
$ws.Range("N5").Value = -5171
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()
$ws.Range("H81").Value = 11373848
$ws.Range("I81").Value = 15158464
$ws.Range("J81").Value = 19999.5
$ws.Range("K81").Value = 30316928
$ws.Range("L81").Value = 39999
$ws.Range("M81").Value = -30315867
$ws.Range("N81").Value = -42121
$ws.Range("H84").Value = 11373848
$ws.Range("I84").Value = 15158464
$ws.Range("J84").Value = 19999.5
$ws.Range("K84").Value = 151584640
$ws.Range("L84").Value = 199995
$ws.Range("M84").Value = -151579336
$ws.Range("N84").Value = -210603
$ws.Range("H107").Value = 385.14285
$ws.Range("I107").Value = 367.875
$ws.Range("K107").Value = 1103.625
$ws.Range("M107").Value = 816.375
$ws.Range("H113").Value = 661.125
$ws.Range("I113").Value = 397.8
$ws.Range("K113").Value = 1193.4
$ws.Range("M113").Value = 976.5999999999999
$ws.Range("H122").Value = 6951.324
$ws.Range("I122").Value = 5215.4585
$ws.Range("J122").Value = 10156
$ws.Range("K122").Value = 15646.3755
$ws.Range("L122").Value = 30468
$ws.Range("M122").Value = -13196.3755
$ws.Range("N122").Value = -35368
$ws.Range("H126").Value = 11736
$ws.Range("J126").Value = 13997.333
$ws.Range("L126").Value = 41991.999
$ws.Range("N126").Value = -46931.999
$ws.Range("H132").Value = 2977.842
$ws.Range("I132").Value = 3031.8
$ws.Range("K132").Value = 9095.400000000001
$ws.Range("M132").Value = -6565.400000000001
$ws.Range("H136").Value = 2184.923
$ws.Range("I136").Value = 2127.9092
$ws.Range("K136").Value = 6383.7276
$ws.Range("M136").Value = -3833.7276
